# Append the two new cost-log rows written by this build run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A77").Value = "2023-12-07 17:15:07"
$ws.Range("B77").Value = 0.0026

$ws.Range("A78").Value = "2023-12-07 17:15:18"
$ws.Range("B78").Value = 0.0004
